$wb = $excel.ActiveWorkbook

# --- Re-arrange / rename sheets -------------------------------------------------
# Fix the "Operazoni" -> "Operazioni" typo (renaming also updates every formula /
# defined name that referenced the old sheet name).
$wb.Worksheets.Item("Operazoni").Name = "Operazioni"

# Add the new placeholder sheets at their target positions.
$wb.Worksheets.Add($wb.Worksheets.Item("Operazioni")).Name = "Sheet_Nuovo2"
$wb.Worksheets.Add($wb.Worksheets.Item("Statistiche")).Name = "Sheet_Nuovo1"
$wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count)).Name = "Sheet_Nuovo"

# Make "Operazioni" the active sheet/tab again.
$wb.Worksheets.Item("Operazioni").Activate()

# --- Content edits on "Operazioni" ----------------------------------------------
$ws = $wb.Worksheets.Item("Operazioni")
$ws.Range("A6").Value = "STOCAZZO"
$ws.Range("A19").Select()

# --- Column width tweak on "Statistiche" ----------------------------------------
# Raw stored column width of 13 chars == ColumnWidth 13 - 0.8333(33) in the
# object model for this workbook's default font (column stays hidden).
$stat = $wb.Worksheets.Item("Statistiche")
$stat.Columns.Item(1).ColumnWidth = 12.166666666666666
